$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Classifieur"
$ws.Range("B1").Value = "nbepochs"
$ws.Range("C1").Value = "batchsize"
$ws.Range("D1").Value = "moteur"
$ws.Range("E1").Value = "base_de_ref"
$ws.Range("F1").Value = "AUC"

# Data rows
$data = @(
    @(1, 10, 4, 2, "sans_Homsap", 0.81280640000000004),
    @(1, 25, 4, 2, "sans_Homsap", 0.89249849999999997),
    @(1, 50, 4, 2, "sans_Homsap", 0.88675090000000001),
    @(1, 75, 4, 2, "sans_Homsap", 0.87661239999999996),
    @(1, 100, 4, 2, "sans_Homsap", 0.89776999999999996),
    @(2, 10, 8, 2, "sans_Homsap", 0.86260650000000005),
    @(2, 25, 8, 2, "sans_Homsap", 0.9019587),
    @(2, 50, 8, 2, "sans_Homsap", 0.89508069999999995),
    @(2, 75, 8, 2, "sans_Homsap", 0.91835639999999996),
    @(2, 100, 8, 2, "sans_Homsap", 0.8850017)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Apply custom font style (Lucida Console, size 8, vertical centered) to AUC column values
$aucRange = $ws.Range("F2:F11")
$aucRange.Font.Name = "Lucida Console"
$aucRange.Font.Size = 8
$aucRange.VerticalAlignment = -4108  # xlCenter

# Column width for column E (base_de_ref)
$ws.Columns.Item(5).ColumnWidth = 12.21875

# Selection matching the saved view state
$ws.Range("G10").Select()
